$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) from 45636 to 45637
# for data rows 2 through 34.
$ws.Range("C2:C34").Value = 45637
